{"js": "await context.sync();\n", "ps1": "$d = $word.ActiveDocument\n"}
